# "first server side push"
# Populates: per-student GPA + self-introduction (学生), refreshed mentor
# roster with new IDs / skill-ratio column + two new mentors (导师); the
# 课题 sheet is untouched content-wise.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "学生" (Students): add gpa (E) + self-introduction (F) columns
# ---------------------------------------------------------------
$wsStudents = $wb.Worksheets.Item("学生")

$wsStudents.Range("E1").Value = "gpa"
$wsStudents.Range("F1").Value = "自我介绍"

$gpa = @{
    2 = 3.21;  3 = 3.3;   4 = 2.5;   5 = 2.6;   6 = 3.4;
    7 = 1.92;  8 = 3.8;   9 = 2.3;   10 = 2.7;  11 = 3.4;
    12 = 3.43; 13 = 3.33; 14 = 3.32; 15 = 3.35; 16 = 3.36;
    17 = 2.01; 18 = 2.01; 19 = 2.03; 20 = 2.04
}

$intro = @{}
$intro[2] = "我是王明"
for ($r = 3; $r -le 20; $r++) { $intro[$r] = "我是李红" }

for ($r = 2; $r -le 20; $r++) {
    $cell = $wsStudents.Range("E$r")
    $cell.NumberFormat = "0.00"
    $cell.Value = $gpa[$r]
    $wsStudents.Range("F$r").Value = $intro[$r]
}

$wsStudents.Columns.Item(5).ColumnWidth = 10.910714285714286
$wsStudents.Range("G16").Select()

# ---------------------------------------------------------------
# Sheet "导师" (Mentors): re-issued ids, wider direction tags, new
# "技艺占比" column, and two freshly added mentors.
# ---------------------------------------------------------------
$wsMentors = $wb.Worksheets.Item("导师")

$wsMentors.Range("D1").Value = "姓别"
$wsMentors.Range("F1").Value = "技艺占比"

$mentorRows = @(
    @{ Row = 2;  Id = 2030513401; Name = "陈伟";   Gender = "男"; Direction = "数据库";                       Ratio = 0.8  },
    @{ Row = 3;  Id = 2030513402; Name = "张得天"; Gender = "男"; Direction = "计算机网络";                   Ratio = 0.9  },
    @{ Row = 4;  Id = 2030513403; Name = "狄岚";   Gender = "女"; Direction = "图形图像处理，程序语言";       Ratio = 0.6  },
    @{ Row = 5;  Id = 2030513404; Name = "张军";   Gender = "男"; Direction = "图形图像处理，程序语言，数据可视化"; Ratio = 0.6  },
    @{ Row = 6;  Id = 2030513405; Name = "陈飞";   Gender = "男"; Direction = "云计算,其它,web";              Ratio = 0.8  },
    @{ Row = 7;  Id = 2030513406; Name = "王士同"; Gender = "男"; Direction = "人工智能,信息安全";            Ratio = 0.95 },
    @{ Row = 8;  Id = 2030513407; Name = "陈秀宏"; Gender = "男"; Direction = "图形图像处理";                 Ratio = 0.85 },
    @{ Row = 9;  Id = 2030513408; Name = "夏鸿斌"; Gender = "男"; Direction = "多媒体,其它";                  Ratio = 0.6  },
    @{ Row = 10; Id = 2030513409; Name = "晏涛";   Gender = "男"; Direction = "图形图像处理，人工智能";       Ratio = 0.7  },
    @{ Row = 11; Id = 2030513410; Name = "赵燕";   Gender = "女"; Direction = "信息可视化，其它";             Ratio = 0.3  },
    @{ Row = 12; Id = 2030513411; Name = "钱鹏江"; Gender = "男"; Direction = "软件工程";                     Ratio = 0.8  },
    @{ Row = 13; Id = 2030513412; Name = "盛新怡"; Gender = "女"; Direction = "人机交互";                     Ratio = 0.4  },
    @{ Row = 14; Id = 2030513413; Name = "律睿敏"; Gender = "男"; Direction = "互动媒体";                     Ratio = 0.5  }
)

foreach ($m in $mentorRows) {
    $r = $m.Row
    $wsMentors.Range("A$r").Value = $m.Id
    $wsMentors.Range("B$r").Value = $m.Name
    $wsMentors.Range("C$r").Value = $m.Id
    $wsMentors.Range("D$r").Value = $m.Gender
    $wsMentors.Range("E$r").Value = $m.Direction
    $wsMentors.Range("F$r").Value = $m.Ratio
}

$wsMentors.Columns.Item(5).ColumnWidth = 32.625
$wsMentors.Activate()
$wsMentors.Range("C17").Select()
